$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Remove the "PPT / 3959222 / KAREN ... / 2211" row (old row 20). This shifts
# rows 21+ up by one, so the last data row (with its distinct bottom-border
# style) becomes row 20, and the closing signature block moves from rows
# 26-27 to rows 25-26.
$ws.Rows.Item(20).Delete()

# Update the 5 remaining worker rows (16-20) with the refreshed data.
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1235045585"
$ws.Range("D16").Value = "JOSE MANUEL CASTRO MONCARIS"
$ws.Range("E16").Value = "2210"
$ws.Range("F16").Value = 40000
$ws.Range("G16").Value = 1000000

$ws.Range("B17").Value = "PPT"
$ws.Range("C17").Value = "3959222"
$ws.Range("D17").Value = "KAREN AURIMAR VALERA MONTERO"
$ws.Range("E17").Value = "2210"
$ws.Range("F17").Value = 40000
$ws.Range("G17").Value = 1533560

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1235045585"
$ws.Range("D18").Value = "JOSE MANUEL CASTRO MONCARIS"
$ws.Range("E18").Value = "2211"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 1000000

$ws.Range("B19").Value = "PPT"
$ws.Range("C19").Value = "3959222"
$ws.Range("D19").Value = "KAREN AURIMAR VALERA MONTERO"
$ws.Range("E19").Value = "2211"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1533560

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143263187"
$ws.Range("D20").Value = "JOSE ALBERTO LAMO CASTAÑEDA"
$ws.Range("E20").Value = "2305"
$ws.Range("F20").Value = 3712
$ws.Range("G20").Value = 1392000

# Refresh the account-summary figures for the new worker set.
$ws.Range("E11").Value = 163712
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 3
